$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> column letter -> new value
$changes = @{
    2 = @{ "D" = "23.414.84"; "E" = "  -0.27%  " }
    3 = @{ "D" = "1.635.49"; "E" = "  -0.67%  " }
    4 = @{ "D" = "1.001"; "E" = "  +0.17%  " }
    5 = @{ "D" = "1.000"; "E" = "  +0.14%  " }
    6 = @{ "D" = "303.45"; "E" = "  -0.39%  " }
    7 = @{ "D" = "0.3774"; "E" = "  +0.63%  " }
    8 = @{ "D" = "51.88"; "E" = "  -0.20%  " }
    9 = @{ "D" = "0.3624"; "E" = "  -0.20%  " }
    10 = @{ "D" = "0.08164"; "E" = "  +0.59%  " }
    11 = @{ "D" = "1.226"; "E" = "  -1.87%  " }
    12 = @{ "D" = "1.001"; "E" = "  +0.17%  " }
    13 = @{ "D" = "22.33"; "E" = "  -2.41%  " }
    14 = @{ "D" = "6.473"; "E" = "  -2.44%  " }
    15 = @{ "D" = "7.346"; "E" = "  +0.58%  " }
    16 = @{ "D" = "0.00001241"; "E" = "  -2.23%  " }
    17 = @{ "D" = "1.631.04"; "E" = "  -0.21%  " }
    18 = @{ "D" = "94.88"; "E" = "  +0.15%  " }
    19 = @{ "D" = "0.06927"; "E" = "  +0.86%  " }
    20 = @{ "D" = "17.58"; "E" = "  -3.46%  " }
    21 = @{ "D" = "6.569"; "E" = "  +0.39%  " }
    22 = @{ "D" = "1.000"; "E" = "  +0.12%  " }
    23 = @{ "D" = "12.52"; "E" = "  -2.32%  " }
    24 = @{ "D" = "23.424.01"; "E" = "  -0.22%  " }
    25 = @{ "D" = "2.497"; "E" = "  +3.72%  " }
    26 = @{ "D" = "3.077"; "E" = "  -2.35%  " }
    27 = @{ "D" = "21.14"; "E" = "  -0.34%  " }
    28 = @{ "D" = "151.02"; "E" = "  +0.18%  " }
    29 = @{ "D" = "5.277"; "E" = "  -0.44%  " }
    30 = @{ "D" = "133.07"; "E" = "  -2.00%  " }
    31 = @{ "D" = "1.812.86"; "E" = "  -0.20%  " }
    32 = @{ "D" = "6.643"; "E" = "  -2.69%  " }
    33 = @{ "D" = "2.156"; "E" = "  -5.69%  " }
    34 = @{ "D" = "1.050"; "E" = "  +9.91%  " }
    35 = @{ "D" = "11.45"; "E" = "  +8.59%  " }
    36 = @{ "D" = "0.02762"; "E" = "  -1.76%  " }
    37 = @{ "B" = "Stellar"; "C" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; "D" = "0.08775"; "E" = "  -0.64%  " }
    38 = @{ "B" = "Algorand"; "C" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; "D" = "0.2489"; "E" = "  -1.79%  " }
    39 = @{ "D" = "0.07110"; "E" = "  -2.32%  " }
    40 = @{ "D" = "6.002"; "E" = "  -3.98%  " }
    41 = @{ "B" = "TrustWalletToken"; "C" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; "D" = "1.344"; "E" = "  -1.89%  " }
    42 = @{ "B" = "TheSandbox"; "C" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; "D" = "0.6999"; "E" = "  -1.02%  " }
    43 = @{ "D" = "15.94"; "E" = "  -1.64%  " }
    44 = @{ "D" = "12.11"; "E" = "  -2.99%  " }
    45 = @{ "D" = "0.6482"; "E" = "  -0.92%  " }
    46 = @{ "D" = "0.9999"; "E" = "  +0.13%  " }
    47 = @{ "D" = "2.273"; "E" = "  -2.71%  " }
    48 = @{ "D" = "3.971"; "E" = "  -0.86%  " }
    49 = @{ "D" = "0.07974"; "E" = "  -0.24%  " }
    50 = @{ "D" = "126.61"; "E" = "  -2.00%  " }
    51 = @{ "D" = "1.185"; "E" = "  -1.83%  " }
}

foreach ($row in $changes.Keys) {
    $rowChanges = $changes[$row]
    foreach ($col in $rowChanges.Keys) {
        $cell = $ws.Range("$col$row")
        if ($col -eq "D") {
            # Force text storage so numeric-looking strings (e.g. "1.001", "303.45")
            # are not coerced into numbers by Excel, matching the original inlineStr text cells.
            $cell.NumberFormat = "@"
            $cell.Value = $rowChanges[$col]
            $cell.ClearFormats()
        } else {
            $cell.Value = $rowChanges[$col]
        }
    }
}
